# Fill the "Tema" column (K) with "Meio Ambiente" for rows 2 through 203,
# matching rows that previously had an empty inline string in that cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2:K203").Value = "Meio Ambiente"
